$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 hold two distinct observation records that were
# mis-assigned; this swaps the per-record data (id, count, coordinates,
# start/end dates, and the public comment) between the two rows while
# leaving every shared attribute (species, location names, etc.) in place.
#
# Using Range.Copy (rather than Range.Value = ...) preserves each cell's
# original data type (numeric vs. text) instead of letting Excel's normal
# literal-entry auto-detection re-interpret numeric-looking or date-looking
# text as numbers/dates.

$scratch = $ws.Range("ZZ1")

$swapCols = @("A", "I", "Q", "R", "Y", "AA")
foreach ($col in $swapCols) {
    $cell5 = $ws.Range($col + "5")
    $cell6 = $ws.Range($col + "6")
    $cell5.Copy($scratch)
    $cell6.Copy($cell5)
    $scratch.Copy($cell6)
    $scratch.Clear()
}

# AC6 ("Riklig") belongs with the record that moved into row 5; row 6's
# record had no public comment.
$ws.Range("AC5").Value = "Riklig"
$ws.Range("AC6").Clear()
